$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 118.833336
$ws.Range("I2").Value = 86.59999999999999
$ws.Range("K2").Value = 86.59999999999999
$ws.Range("M2").Value = 26.40000000000001

$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H80").Value = 2123.2222
$ws.Range("I80").Value = 602.0833
$ws.Range("K80").Value = 1806.2499
$ws.Range("M80").Value = -808.2499

$ws.Range("H83").Value = 2123.2222
$ws.Range("I83").Value = 602.0833
$ws.Range("K83").Value = 5418.7497
$ws.Range("M83").Value = -426.7497000000003

$ws.Range("H92").Value = 1888.0834
$ws.Range("I92").Value = 497.25
$ws.Range("K92").Value = 497.25
$ws.Range("M92").Value = 750.75

$ws.Range("J96").Value = 16000
$ws.Range("L96").Value = 48000
$ws.Range("N96").Value = -50746

$ws.Range("H100").Value = 3868.9167
$ws.Range("I100").Value = 2800
$ws.Range("J100").Value = 4082.7
$ws.Range("K100").Value = 2800
$ws.Range("L100").Value = 4082.7
$ws.Range("M100").Value = -2259
$ws.Range("N100").Value = -5164.7

$ws.Range("H113").Value = 15000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 15000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 15000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -21508

$ws.Range("H138").Value = 2487.8262
$ws.Range("I138").Value = 1841.72
$ws.Range("J138").Value = 3257
$ws.Range("K138").Value = 5525.16
$ws.Range("L138").Value = 9771
$ws.Range("M138").Value = -385.1599999999999
$ws.Range("N138").Value = -20051

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4280.34
$ws.Range("I32").Value = 3518.9148
$ws.Range("J32").Value = 16209.333
$ws.Range("K32").Value = 3518.9148
$ws.Range("L32").Value = 16209.333
$ws.Range("M32").Value = -3231.9148
$ws.Range("N32").Value = -16783.333

$ws.Range("H61").Value = 8545.091
$ws.Range("I61").Value = 5599.5713
$ws.Range("K61").Value = 5599.5713
$ws.Range("M61").Value = -5387.5713

$ws.Range("H110").Value = 3856.3333
$ws.Range("J110").Value = 9880.5
$ws.Range("L110").Value = 9880.5
$ws.Range("N110").Value = -13970.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 8545.091
$ws.Range("I136").Value = 5599.5713
$ws.Range("K136").Value = 16798.7139
$ws.Range("M136").Value = -14248.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7814.1904
$ws.Range("I105").Value = 3715.6667
$ws.Range("K105").Value = 3715.6667
$ws.Range("M105").Value = -1968.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32512.945
$ws.Range("I31").Value = 3236.652
$ws.Range("J31").Value = 84309.46000000001
$ws.Range("K31").Value = 3236.652
$ws.Range("L31").Value = 84309.46000000001
$ws.Range("M31").Value = -2941.652
$ws.Range("N31").Value = -84899.46000000001

$ws.Range("H34").Value = 32512.945
$ws.Range("I34").Value = 3236.652
$ws.Range("J34").Value = 84309.46000000001
$ws.Range("K34").Value = 3236.652
$ws.Range("L34").Value = 84309.46000000001
$ws.Range("M34").Value = -3034.652
$ws.Range("N34").Value = -84713.46000000001

$ws.Range("H43").Value = 7000
$ws.Range("J43").Value = 7000
$ws.Range("L43").Value = 7000
$ws.Range("N43").Value = -7368

$ws.Range("H74").Value = 160020.6
$ws.Range("J74").Value = 160020.6
$ws.Range("L74").Value = 160020.6
$ws.Range("N74").Value = -161768.6

$ws.Range("H77").Value = 160020.6
$ws.Range("J77").Value = 160020.6
$ws.Range("L77").Value = 480061.8
$ws.Range("N77").Value = -488797.8

$ws.Range("H99").Value = 2114.25
$ws.Range("I99").Value = 2150
$ws.Range("J99").Value = 2007
$ws.Range("K99").Value = 2150
$ws.Range("L99").Value = 2007
$ws.Range("M99").Value = -652
$ws.Range("N99").Value = -5003

$ws.Range("H100").Value = 59999.668
$ws.Range("J100").Value = 59999.668
$ws.Range("L100").Value = 59999.668
$ws.Range("N100").Value = -62163.668

$ws.Range("H101").Value = 7000
$ws.Range("J101").Value = 7000
$ws.Range("L101").Value = 7000
$ws.Range("N101").Value = -13490

$ws.Range("H126").Value = 2114.25
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 2007
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 6021
$ws.Range("M126").Value = -3980
$ws.Range("N126").Value = -10961

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 9091.091
$ws.Range("J69").Value = 9999
$ws.Range("L69").Value = 29997
$ws.Range("N69").Value = -31619

$ws.Range("H72").Value = 9091.091
$ws.Range("J72").Value = 9999
$ws.Range("L72").Value = 89991
$ws.Range("N72").Value = -98103

$ws.Range("H93").Value = 12013.5
$ws.Range("I93").Value = 5000
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = -13128

$ws.Range("H139").Value = 2692.625
$ws.Range("I139").Value = 1603.2667
$ws.Range("J139").Value = 19033
$ws.Range("K139").Value = 4809.800099999999
$ws.Range("L139").Value = 57099
$ws.Range("M139").Value = 330.1999000000005
$ws.Range("N139").Value = -67379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3431.5715
$ws.Range("I122").Value = 2796.875
$ws.Range("K122").Value = 8390.625
$ws.Range("M122").Value = -5940.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6454.091
$ws.Range("I7").Value = 2763.7778
$ws.Range("J7").Value = 9008.923000000001
$ws.Range("K7").Value = 2763.7778
$ws.Range("L7").Value = 9008.923000000001
$ws.Range("M7").Value = -2651.7778
$ws.Range("N7").Value = -9232.923000000001

$ws.Range("H61").Value = 6950.6
$ws.Range("I61").Value = 1401
$ws.Range("J61").Value = 10650.333
$ws.Range("K61").Value = 1401
$ws.Range("L61").Value = 10650.333
$ws.Range("M61").Value = -1199
$ws.Range("N61").Value = -11054.333

$ws.Range("H68").Value = 3636.5
$ws.Range("I68").Value = 2968.3914
$ws.Range("K68").Value = 2968.3914
$ws.Range("M68").Value = -2219.3914

$ws.Range("H71").Value = 3636.5
$ws.Range("I71").Value = 2968.3914
$ws.Range("K71").Value = 14841.957
$ws.Range("M71").Value = -11097.957

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H113").Value = 6950.6
$ws.Range("I113").Value = 1401
$ws.Range("J113").Value = 10650.333
$ws.Range("K113").Value = 1401
$ws.Range("L113").Value = 10650.333
$ws.Range("M113").Value = 769
$ws.Range("N113").Value = -14990.333

$ws.Range("H126").Value = 6454.091
$ws.Range("I126").Value = 2763.7778
$ws.Range("J126").Value = 9008.923000000001
$ws.Range("K126").Value = 8291.3334
$ws.Range("L126").Value = 27026.769
$ws.Range("M126").Value = -5821.3334
$ws.Range("N126").Value = -31966.769

$ws.Range("H136").Value = 9564.462
$ws.Range("I136").Value = 4980
$ws.Range("K136").Value = 14940
$ws.Range("M136").Value = -12390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8862.5
$ws.Range("J122").Value = 16490.637
$ws.Range("L122").Value = 49471.91099999999
$ws.Range("N122").Value = -54371.91099999999

$ws.Range("H136").Value = 3244.7778
$ws.Range("I136").Value = 2317.7058
$ws.Range("K136").Value = 6953.117400000001
$ws.Range("M136").Value = -4403.117400000001
